# "element buff skill fin" - rebalance the Buff sheet's "Mark" numeric
# column (L, rows 4-23) by dividing every data-row value by 20 (tuning
# the elemental-buff stack thresholds down to the new, smaller scale),
# then carry over the residual view/formatting state (active selection
# and a best-fit column width hanging off an empty column) that Excel
# re-wrote when the sheet was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Buff")

# --- L4:L23 numeric rebalance (divide by 20) ---
$values = @{
    4  = 25
    5  = 30
    6  = 30
    7  = 25
    8  = 50
    9  = 20
    10 = 100
    11 = 10
    12 = 30
    13 = 125
    14 = 25
    15 = 50
    16 = 75
    17 = 75
    18 = 30
    19 = 20
    20 = 40
    21 = 20
    22 = 10
    23 = 10
}

foreach ($row in $values.Keys) {
    $ws.Range("L$row").Value = $values[$row]
}

# --- stray "bestFit" column width metadata shifts from column P (16) to
# column O (15); both are empty, so deleting the (empty) column O shifts
# column P's width record left by one, landing exactly on column O. ---
$ws.Columns.Item(15).Delete()

# --- restore the last active selection recorded in the saved view state ---
$ws.Activate()
$ws.Range("K12").Select()
